# Update the "想去人数" (want-to-go count) figures that changed between
# the previous gh-pages data pull and the new one (commit 456a3b4).
# Both the "展览" sheet and its "全部类型" mirror need the same updates.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F5"  = 4688
    "F7"  = 404
    "F8"  = 1412
    "F11" = 1219
    "F13" = 678
    "F15" = 60
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
